$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 69: Grabherr, Gottfried & Pauli (2010) ----
$ws.Range("A69").Value = "Climate Change Impacts in Alpine Environments"
$ws.Range("B69").Value = "Geography compass"
$ws.Range("C69").Value = 2010
$ws.Range("D69").Value = "Grabherr, Gottfried & Pauli"
$ws.Range("F69").Value = "≈"
$ws.Range("G69").Value = "Climate change in Apline environment"
$ws.Range("I69").Value = "yes"
$ws.Range("J69").Value = "yes"
$ws.Range("K69").Value = "yes"

# ---- Row 70: Kristensen & al. (2016) - TMB package ----
$ws.Range("A70").Value = "TMB: Automatic Differentiation and Laplace Approximation"
$ws.Range("B70").Value = "Journal of Statistical Software"
$ws.Range("C70").Value = 2016
$ws.Range("D70").Value = "Kristensen & al."
$ws.Range("E70").Value = "Software/Package"
$ws.Range("F70").Value = "NA"
$ws.Range("G70").Value = "TMB package"
$ws.Range("I70").Value = "yes"
$ws.Range("J70").Value = "yes"
$ws.Range("K70").Value = "yes"

# ---- Fix up cell formatting to mirror the rest of the table ----
# F69 picks up the "approx" style used elsewhere in column F (e.g. F45)
$ws.Range("F45").Copy() | Out-Null
$ws.Range("F69").PasteSpecial(-4122) | Out-Null

# E70 picks up the "Software/Package" style used elsewhere in column E (e.g. E67)
$ws.Range("E67").Copy() | Out-Null
$ws.Range("E70").PasteSpecial(-4122) | Out-Null

# F70 picks up the "NA" style used elsewhere in column F (e.g. F68)
$ws.Range("F68").Copy() | Out-Null
$ws.Range("F70").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
